$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 60).Value = "26/04/2020"
$ws.Cells.Item(2, 16).Value = $null
$ws.Cells.Item(2, 52).Value = $null
$ws.Cells.Item(2, 54).Value = 1
$ws.Cells.Item(2, 60).Value = $null
$ws.Cells.Item(3, 54).Value = $null
$ws.Cells.Item(3, 55).Value = 5
$ws.Cells.Item(3, 56).Value = 3
$ws.Cells.Item(3, 58).Value = $null
$ws.Cells.Item(3, 59).Value = 1
$ws.Cells.Item(3, 60).Value = $null
$ws.Cells.Item(4, 54).Value = 81
$ws.Cells.Item(4, 55).Value = 2
$ws.Cells.Item(4, 56).Value = 6
$ws.Cells.Item(4, 57).Value = 2
$ws.Cells.Item(4, 58).Value = $null
$ws.Cells.Item(4, 59).Value = 2
$ws.Cells.Item(4, 60).Value = 21
$ws.Cells.Item(5, 60).Value = $null
$ws.Cells.Item(6, 59).Value = $null
$ws.Cells.Item(6, 60).Value = 1
$ws.Cells.Item(7, 54).Value = $null
$ws.Cells.Item(7, 56).Value = 2
$ws.Cells.Item(7, 57).Value = 4
$ws.Cells.Item(7, 59).Value = $null
$ws.Cells.Item(7, 60).Value = 6
$ws.Cells.Item(8, 54).Value = $null
$ws.Cells.Item(8, 55).Value = 1
$ws.Cells.Item(8, 58).Value = $null
$ws.Cells.Item(8, 59).Value = 5
$ws.Cells.Item(8, 60).Value = 2
$ws.Cells.Item(9, 55).Value = $null
$ws.Cells.Item(9, 56).Value = 1
$ws.Cells.Item(9, 58).Value = $null
$ws.Cells.Item(9, 60).Value = 1
$ws.Cells.Item(10, 54).Value = $null
$ws.Cells.Item(10, 55).Value = 1
$ws.Cells.Item(10, 56).Value = 9
$ws.Cells.Item(10, 57).Value = 23
$ws.Cells.Item(10, 58).Value = $null
$ws.Cells.Item(10, 59).Value = 25
$ws.Cells.Item(10, 60).Value = 23
$ws.Cells.Item(11, 56).Value = $null
$ws.Cells.Item(11, 57).Value = 2
$ws.Cells.Item(11, 60).Value = $null
$ws.Cells.Item(12, 54).Value = $null
$ws.Cells.Item(12, 55).Value = 3
$ws.Cells.Item(12, 57).Value = 8
$ws.Cells.Item(12, 58).Value = $null
$ws.Cells.Item(12, 59).Value = 7
$ws.Cells.Item(12, 60).Value = 4
$ws.Cells.Item(13, 54).Value = $null
$ws.Cells.Item(13, 55).Value = 1
$ws.Cells.Item(13, 56).Value = 2
$ws.Cells.Item(13, 58).Value = $null
$ws.Cells.Item(13, 59).Value = 1
$ws.Cells.Item(13, 60).Value = 7
$ws.Cells.Item(14, 60).Value = $null
$ws.Cells.Item(15, 54).Value = $null
$ws.Cells.Item(15, 55).Value = 6
$ws.Cells.Item(15, 56).Value = 13
$ws.Cells.Item(15, 57).Value = 5
$ws.Cells.Item(15, 58).Value = $null
$ws.Cells.Item(15, 59).Value = 11
$ws.Cells.Item(15, 60).Value = 6
$ws.Cells.Item(16, 59).Value = $null
$ws.Cells.Item(16, 60).Value = 1
$ws.Cells.Item(17, 54).Value = 1
$ws.Cells.Item(17, 55).Value = 9
$ws.Cells.Item(17, 56).Value = 12
$ws.Cells.Item(17, 57).Value = 7
$ws.Cells.Item(17, 58).Value = $null
$ws.Cells.Item(17, 59).Value = 5
$ws.Cells.Item(17, 60).Value = 13
$ws.Cells.Item(18, 16).Value = 9
$ws.Cells.Item(18, 54).Value = $null
$ws.Cells.Item(18, 55).Value = 39
$ws.Cells.Item(18, 56).Value = 6
$ws.Cells.Item(18, 57).Value = 20
$ws.Cells.Item(18, 58).Value = 1
$ws.Cells.Item(18, 59).Value = 32
$ws.Cells.Item(18, 60).Value = 7
$ws.Cells.Item(19, 60).Value = $null
$ws.Cells.Item(20, 54).Value = $null
$ws.Cells.Item(20, 55).Value = 9
$ws.Cells.Item(20, 56).Value = 10
$ws.Cells.Item(20, 58).Value = $null
$ws.Cells.Item(20, 59).Value = 8
$ws.Cells.Item(20, 60).Value = 19
$ws.Cells.Item(21, 54).Value = $null
$ws.Cells.Item(21, 56).Value = 28
$ws.Cells.Item(21, 57).Value = 22
$ws.Cells.Item(21, 58).Value = $null
$ws.Cells.Item(21, 59).Value = 17
$ws.Cells.Item(21, 60).Value = 82
$ws.Cells.Item(22, 60).Value = $null
$ws.Cells.Item(23, 60).Value = $null
$ws.Cells.Item(24, 60).Value = $null
$ws.Cells.Item(25, 54).Value = 3
$ws.Cells.Item(25, 59).Value = $null
$ws.Cells.Item(25, 60).Value = 1
$ws.Cells.Item(26, 60).Value = $null
$ws.Cells.Item(27, 54).Value = $null
$ws.Cells.Item(27, 55).Value = 5
$ws.Cells.Item(27, 58).Value = $null
$ws.Cells.Item(27, 59).Value = 3
$ws.Cells.Item(27, 60).Value = 1
$ws.Cells.Item(28, 54).Value = 58
$ws.Cells.Item(28, 55).Value = 7
$ws.Cells.Item(28, 56).Value = 4
$ws.Cells.Item(28, 57).Value = 5
$ws.Cells.Item(28, 58).Value = $null
$ws.Cells.Item(28, 59).Value = 20
$ws.Cells.Item(28, 60).Value = 14
$ws.Cells.Item(29, 54).Value = $null
$ws.Cells.Item(29, 55).Value = 9
$ws.Cells.Item(29, 56).Value = 4
$ws.Cells.Item(29, 57).Value = 8
$ws.Cells.Item(29, 58).Value = $null
$ws.Cells.Item(29, 59).Value = 17
$ws.Cells.Item(29, 60).Value = 57
$ws.Cells.Item(30, 54).Value = $null
$ws.Cells.Item(30, 55).Value = 14
$ws.Cells.Item(30, 56).Value = 8
$ws.Cells.Item(30, 57).Value = 3
$ws.Cells.Item(30, 58).Value = $null
$ws.Cells.Item(30, 59).Value = 7
$ws.Cells.Item(30, 60).Value = 20
$ws.Cells.Item(31, 60).Value = $null
$ws.Cells.Item(32, 54).Value = $null
$ws.Cells.Item(32, 55).Value = 7
$ws.Cells.Item(32, 56).Value = 16
$ws.Cells.Item(32, 57).Value = 7
$ws.Cells.Item(32, 58).Value = $null
$ws.Cells.Item(32, 59).Value = 24
$ws.Cells.Item(32, 60).Value = 8
$ws.Cells.Item(33, 54).Value = 2
$ws.Cells.Item(33, 55).Value = $null
$ws.Cells.Item(33, 57).Value = 1
$ws.Cells.Item(33, 60).Value = $null
$ws.Cells.Item(34, 54).Value = 40
$ws.Cells.Item(34, 55).Value = 5
$ws.Cells.Item(34, 57).Value = 3
$ws.Cells.Item(34, 58).Value = $null
$ws.Cells.Item(34, 59).Value = 1
$ws.Cells.Item(34, 60).Value = 15
